# Regenerate merged AHB files
# 1. Rename the shared-string column headers from the *_old/*_new suffix
#    convention to the *_FV2310/*_FV2404 convention.
# 2. Freeze the header row (row 1).
# 3. Turn the A1:U81 range into an Excel Table ("Table1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$renames = @{
    "Segmentname_old"          = "Segmentname_FV2310"
    "Segmentgruppe_old"        = "Segmentgruppe_FV2310"
    "Segment_old"               = "Segment_FV2310"
    "Datenelement_old"          = "Datenelement_FV2310"
    "Segment ID_old"            = "Segment ID_FV2310"
    "Code_old"                  = "Code_FV2310"
    "Qualifier_old"             = "Qualifier_FV2310"
    "Beschreibung_old"          = "Beschreibung_FV2310"
    "Bedingungsausdruck_old"    = "Bedingungsausdruck_FV2310"
    "Bedingung_old"             = "Bedingung_FV2310"
    "Segmentname_new"           = "Segmentname_FV2404"
    "Segmentgruppe_new"         = "Segmentgruppe_FV2404"
    "Segment_new"                = "Segment_FV2404"
    "Datenelement_new"           = "Datenelement_FV2404"
    "Segment ID_new"             = "Segment ID_FV2404"
    "Code_new"                   = "Code_FV2404"
    "Qualifier_new"              = "Qualifier_FV2404"
    "Beschreibung_new"           = "Beschreibung_FV2404"
    "Bedingungsausdruck_new"     = "Bedingungsausdruck_FV2404"
    "Bedingung_new"              = "Bedingung_FV2404"
}

for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $old = $cell.Text
    if ($renames.ContainsKey($old)) {
        $cell.Value = $renames[$old]
    }
}

# Freeze the header row.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Convert the used range into a table.
$rng = $ws.Range("A1:U81")
$tbl = $ws.ListObjects.Add(1, $rng, 0, 1)
$tbl.Name = "Table1"
